$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin price / 1h volume figures (and one Filecoin <-> PancakeSwap
# row-order swap) pulled in by the scheduled GitHub Actions data update.
# Price-column values are free-form text (e.g. "34.401.12", "0.0964") even
# when they happen to look numeric, so each such cell is force-formatted as
# text right before the write -- this keeps e.g. "11.30" from being silently
# re-interpreted as the number 11.3 (dropping the trailing zero).

$ws.Range('D2').Value = '34.401.12'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.804.35'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.37'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.602'
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.31'
$ws.Range('E8').Value = '  +3.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.293'
$ws.Range('E9').Value = '  -2.27%  '
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0964'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').Value = '2.062.80'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.30'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').Value = '1.794.57'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('D16').Value = '34.392.56'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.44'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.55'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.52'
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.22'
$ws.Range('E21').Value = '  -2.37%  '
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('E24').Value = '  +5.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.09'
$ws.Range('E26').Value = '  +3.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.31'
$ws.Range('E27').Value = '  +3.63%  '
$ws.Range('E28').Value = '  +2.18%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.91'
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.23'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.79'
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0516'
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.78'
$ws.Range('E34').Value = '  -2.76%  '
$ws.Range('D35').Value = '1.360.33'
$ws.Range('E35').Value = '  -2.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.652'
$ws.Range('E36').Value = '  -3.91%  '
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('E38').Value = '  -6.29%  '
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.78'
$ws.Range('E41').Value = '  -1.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '80.62'
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('E44').Value = '  +5.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.31'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0497'
$ws.Range('E46').Value = '  -2.84%  '
$ws.Range('D47').Value = '1.965.13'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('E48').Value = '  -3.60%  '
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.33'
$ws.Range('E50').Value = '  -2.08%  '
$ws.Range('D51').Value = '0.0₆0122'
$ws.Range('E51').Value = '  -5.54%  '
